$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to remain text before assigning,
# so Excel does not auto-convert numeric-looking strings to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '33.858.72'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.778.63'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '221.00'
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").Value = '0.550'
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '31.02'
$ws.Range("E8").Value = '  -5.01%  '
$ws.Range("D9").Value = '0.285'
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("D10").Value = '0.0708'
$ws.Range("E10").Value = '  +5.81%  '
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("D12").Value = '2.032.41'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.769.69'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '10.47'
$ws.Range("E14").Value = '  -5.02%  '
$ws.Range("D15").Value = '0.625'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '33.835.93'
$ws.Range("E16").Value = '  -2.27%  '
$ws.Range("D17").Value = '4.21'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '67.92'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").Value = '244.66'
$ws.Range("E19").Value = '  -3.67%  '
$ws.Range("D20").Value = '0.0₃0773'
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '10.63'
$ws.Range("E22").Value = '  +2.08%  '
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").Value = '  -3.72%  '
$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D25").Value = '157.43'
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("D26").Value = '16.38'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("D28").Value = '0.112'
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '0.0520'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = '3.70'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '3.49'
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").Value = '1.80'
$ws.Range("E34").Value = '  -2.56%  '
$ws.Range("D35").Value = '1.394.92'
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("D36").Value = '0.639'
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '0.0185'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").Value = '0.932'
$ws.Range("E39").Value = '  +3.58%  '
$ws.Range("D40").Value = '2.34'
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("D41").Value = '78.97'
$ws.Range("E41").Value = '  -4.94%  '
$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  -3.88%  '
$ws.Range("D43").Value = '2.11'
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("D44").Value = '5.93'
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = '0.0489'
$ws.Range("E45").Value = '  -3.03%  '
$ws.Range("D46").Value = '1.03'
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").Value = '1.930.63'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = '104.78'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").Value = '0.996'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("D50").Value = '11.70'
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").Value = '0.0₆0120'
$ws.Range("E51").Value = '  -1.62%  '

# Remove the temporary text-number-format so the cell style matches the original (no explicit style index).
$ws.Range("D2:E51").ClearFormats()
